$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation for numeric-looking
# strings (so "216.61" etc. stay text, matching the source inline/shared strings,
# instead of being auto-coerced to a float by Excel). Style is restored afterwards
# so no visible formatting changes.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
$ws.Cells.Item(2, 4).Value = '26.027.60'
$ws.Cells.Item(2, 5).Value = '  +0.62%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.644.28'
$ws.Cells.Item(3, 5).Value = '  +0.91%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +1.09%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '216.61'
$ws.Cells.Item(5, 5).Value = '  +1.09%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +1.10%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.90%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.75%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +1.36%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '19.66'
$ws.Cells.Item(10, 5).Value = '  +0.27%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.0795'
$ws.Cells.Item(11, 5).Value = '  +1.04%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(12, 4) '4.30'
$ws.Cells.Item(12, 5).Value = '  +1.62%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '1.872.11'
$ws.Cells.Item(13, 5).Value = '  +0.75%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.651.44'
$ws.Cells.Item(14, 5).Value = '  +1.80%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +0.40%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₃0764'
$ws.Cells.Item(16, 5).Value = '  +1.08%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.81%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '26.029.15'
$ws.Cells.Item(18, 5).Value = '  +0.62%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.92%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '192.96'
$ws.Cells.Item(20, 5).Value = '  +0.23%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.50%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '9.95'
$ws.Cells.Item(22, 5).Value = '  +0.16%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.07%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.35%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Monero'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Cells.Item(25, 4) '144.57'
$ws.Cells.Item(25, 5).Value = '  +1.70%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Cells.Item(26, 4) '0.131'
$ws.Cells.Item(26, 5).Value = '  +4.11%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +1.07%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.63%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '15.54'
$ws.Cells.Item(29, 5).Value = '  +0.58%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.10%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '0.0498'
$ws.Cells.Item(31, 5).Value = '  -0.08%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.32%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.29%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +2.47%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -3.02%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '0.907'
$ws.Cells.Item(36, 5).Value = '  +0.77%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.133.66'
$ws.Cells.Item(37, 5).Value = '  +0.07%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.541'
$ws.Cells.Item(38, 5).Value = '  -1.33%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.61%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.75%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '5.50'
$ws.Cells.Item(41, 5).Value = '  +0.80%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '99.52'
$ws.Cells.Item(42, 5).Value = '  +0.55%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '0.798'
$ws.Cells.Item(43, 5).Value = '  -0.56%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '1.781.63'
$ws.Cells.Item(44, 5).Value = '  +0.77%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.0₆0116'
$ws.Cells.Item(45, 5).Value = '  +3.65%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '56.54'
$ws.Cells.Item(46, 5).Value = '  +0.97%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.0532'
$ws.Cells.Item(47, 5).Value = '  +1.06%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.31%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '7.71'
$ws.Cells.Item(49, 5).Value = '  +1.60%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.83%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.28%  '
